# Applies:
#  1) Slide 5's table ("Google Shape;122;p17") gets a new table style GUID
#     ({62D91FEB-F0E7-412A-9C79-BCF3B7736578} -> {047F71B6-0B34-40D4-8E06-CA2F3937E355}).
#  2) The presentation's theme (ppt/theme/theme1.xml, the one bound to the
#     slide master / presentation) is recolored from the "Integral" / "Red
#     Violet" palette to the stock "Office" palette (same palette that used
#     to live in ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{047F71B6-0B34-40D4-8E06-CA2F3937E355}")

# --- 2) Recolor the theme --------------------------------------------------
function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme color scheme, in ThemeColorScheme.Item(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbLong $officeColors[$i - 1]
}
